$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data (row 18) describing a new test/task item.
$ws.Range("A18").Value = "refactor IMU reading (i.e. methods like read.gyro.z())"
$ws.Range("B18").Value = "N"
$ws.Range("C18").Value = "N"
$ws.Range("D18").Value = "Design"

# Apply the same centered alignment used by the rest of column B/C/D to the new cells.
$ws.Range("B18:D18").HorizontalAlignment = -4108

# Update the active selection like the source workbook (cursor moved to D19 after entry).
$ws.Range("D19").Select()
